$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (3-7); only the header + one data row remain.
$ws.Rows("3:7").Delete()

# --- Header row (row 1): fecha, componentes, cantidad, responsable ---
$ws.Range("A1").Value = "fecha"
$ws.Range("B1").Value = "componentes"
$ws.Range("C1").Value = "cantidad"
# D1 ("responsable") is unchanged.

# --- Data row (row 2): 2024-01-24, c1818, 60, JORGE FUENTES ---
# Column A needs the literal text "2024-01-24" (not an Excel date serial).
# Using Formula with a quoted string forces a text result without Excel's
# auto date-recognition kicking in, then PasteSpecial (values only) bakes
# that text in as a plain value/string cell (no formula, no number format).
$ws.Range("A2").Formula = '="2024-01-24"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B2").Value = "c1818"
$ws.Range("C2").Value = 60
# D2 ("JORGE FUENTES") is unchanged.
